$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts at row 2).
$lastRow = $ws.UsedRange.Rows.Count - 1

# 1) Update the "Förändrad" (changed) date in column C for every data row
#    from 45184 to 45186 (2023-09-15 -> 2023-09-17).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# 2) For rows that contain the link formulas (S, T, V, W, X, Y), add the
#    beteckning (value from column A) as the second, friendly-name
#    argument to the HYPERLINK() function, if it isn't already present.
$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    $id = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($id)) {
        continue
    }

    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $r)
        $f = $cell.Formula
        if ([string]::IsNullOrEmpty($f)) {
            continue
        }
        if ($f.StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $id + '")'
            $cell.Formula = $newFormula
        }
    }
}
